$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $ref, $val)
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$updates = @(
    @{Ref='D2'; Val='19.935.72'},
    @{Ref='E2'; Val='  -7.98%  '},
    @{Ref='D3'; Val='1.402.68'},
    @{Ref='E3'; Val='  -8.51%  '},
    @{Ref='D4'; Val='1.001'},
    @{Ref='E4'; Val='  +0.06%  '},
    @{Ref='E5'; Val='  +0.09%  '},
    @{Ref='D6'; Val='270.04'},
    @{Ref='E6'; Val='  -6.33%  '},
    @{Ref='D7'; Val='0.3656'},
    @{Ref='E7'; Val='  -7.23%  '},
    @{Ref='D8'; Val='0.3028'},
    @{Ref='E8'; Val='  -3.94%  '},
    @{Ref='D9'; Val='39.09'},
    @{Ref='E9'; Val='  -7.14%  '},
    @{Ref='D10'; Val='0.06439'},
    @{Ref='E10'; Val='  -10.00%  '},
    @{Ref='D11'; Val='0.9669'},
    @{Ref='E11'; Val='  -7.34%  '},
    @{Ref='E12'; Val='  +0.11%  '},
    @{Ref='D13'; Val='5.260'},
    @{Ref='E13'; Val='  -6.52%  '},
    @{Ref='D14'; Val='6.061'},
    @{Ref='E15'; Val='  -10.18%  '},
    @{Ref='D16'; Val='1.406.02'},
    @{Ref='E16'; Val='  -8.58%  '},
    @{Ref='D17'; Val='0.000009978'},
    @{Ref='E17'; Val='  -8.53%  '},
    @{Ref='D18'; Val='0.05665'},
    @{Ref='E18'; Val='  -14.05%  '},
    @{Ref='B19'; Val='Dai'},
    @{Ref='C19'; Val='https://coinranking.com/coin/MoTuySvg7+dai-dai'},
    @{Ref='D19'; Val='1.002'},
    @{Ref='E19'; Val='  +0.09%  '},
    @{Ref='B20'; Val='Litecoin'},
    @{Ref='C20'; Val='https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'},
    @{Ref='D20'; Val='71.29'},
    @{Ref='E20'; Val='  -14.02%  '},
    @{Ref='D21'; Val='5.475'},
    @{Ref='E21'; Val='  -10.08%  '},
    @{Ref='D22'; Val='14.14'},
    @{Ref='E22'; Val='  -8.06%  '},
    @{Ref='D23'; Val='10.54'},
    @{Ref='E23'; Val='  -2.47%  '},
    @{Ref='D24'; Val='2.275'},
    @{Ref='E24'; Val='  -4.81%  '},
    @{Ref='D25'; Val='19.923.86'},
    @{Ref='E25'; Val='  -8.01%  '},
    @{Ref='D26'; Val='2.205'},
    @{Ref='E26'; Val='  -5.28%  '},
    @{Ref='D27'; Val='135.52'},
    @{Ref='E27'; Val='  -8.08%  '},
    @{Ref='D28'; Val='16.52'},
    @{Ref='E28'; Val='  -9.60%  '},
    @{Ref='D29'; Val='1.564.35'},
    @{Ref='E29'; Val='  -8.60%  '},
    @{Ref='D30'; Val='107.38'},
    @{Ref='E30'; Val='  -8.14%  '},
    @{Ref='D31'; Val='3.826'},
    @{Ref='E31'; Val='  -20.94%  '},
    @{Ref='D32'; Val='5.195'},
    @{Ref='E32'; Val='  -11.03%  '},
    @{Ref='D33'; Val='0.7924'},
    @{Ref='E33'; Val='  -16.31%  '},
    @{Ref='D34'; Val='0.07571'},
    @{Ref='E34'; Val='  -6.82%  '},
    @{Ref='D35'; Val='8.272'},
    @{Ref='E35'; Val='  -2.64%  '},
    @{Ref='B36'; Val='Frax'},
    @{Ref='C36'; Val='https://coinranking.com/coin/KfWtaeV1W+frax-frax'},
    @{Ref='D36'; Val='1.001'},
    @{Ref='E36'; Val='  +0.02%  '},
    @{Ref='B37'; Val='Hedera'},
    @{Ref='C37'; Val='https://coinranking.com/coin/jad286TjB+hedera-hbar'},
    @{Ref='D37'; Val='0.05698'},
    @{Ref='E37'; Val='  -5.80%  '},
    @{Ref='D38'; Val='4.672'},
    @{Ref='E38'; Val='  -7.93%  '},
    @{Ref='D39'; Val='0.1909'},
    @{Ref='E39'; Val='  -5.18%  '},
    @{Ref='D40'; Val='0.02003'},
    @{Ref='E40'; Val='  -8.87%  '},
    @{Ref='D41'; Val='1.328'},
    @{Ref='E41'; Val='  -7.76%  '},
    @{Ref='E42'; Val='  -6.76%  '},
    @{Ref='D43'; Val='1.046'},
    @{Ref='E43'; Val='  -11.08%  '},
    @{Ref='B44'; Val='PancakeSwap'},
    @{Ref='C44'; Val='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'},
    @{Ref='D44'; Val='3.484'},
    @{Ref='E44'; Val='  -6.37%  '},
    @{Ref='B45'; Val='TheSandbox'},
    @{Ref='C45'; Val='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'},
    @{Ref='D45'; Val='0.5198'},
    @{Ref='E45'; Val='  -9.25%  '},
    @{Ref='D46'; Val='11.92'},
    @{Ref='E46'; Val='  -8.03%  '},
    @{Ref='D47'; Val='0.5009'},
    @{Ref='E47'; Val='  -8.37%  '},
    @{Ref='D48'; Val='108.92'},
    @{Ref='E48'; Val='  -5.49%  '},
    @{Ref='D49'; Val='1.738'},
    @{Ref='E49'; Val='  -6.43%  '},
    @{Ref='D50'; Val='1.003'},
    @{Ref='E50'; Val='  +0.08%  '},
    @{Ref='D51'; Val='1.027'},
    @{Ref='E51'; Val='  -10.83%  '}
)

foreach ($u in $updates) {
    Set-TextValue $ws $u.Ref $u.Val
}
